# Jay Sprout Dev-Resume - "Create HTML skeleton, update ID-index"
#
# Skills section edits:
#   1. JavaScript | HTML | ...            -> JavaScript | D3 | Node | HTML | ...
#   2. ... | Git | XML                    -> ... | Git | XML | SAM | OOP
#   3. "Search Engine Optimization (SEO)" -> "SEO"
#   4. Merge the "SEO | Progressive Web Apps" paragraph with the following
#      "User Experience (UX) | User Interface (UI) | ..." paragraph.
#   5. "User Experience (UX)"             -> " | Agile | Scrum | UX"
#   6. "User Interface (UI) "             -> "UI "
#
$d = $word.ActiveDocument

# 1. Insert "D3 | Node | " right before the first "HTML" (skills line).
$r = $d.Content
$r.Find.Execute("HTML", $true, $false, $false, $false, $false, `
                $true, 1, $false, "", 0)
$ins = $d.Range($r.Start, $r.Start)
$ins.InsertAfter("D3 | Node | ")

# 2. Append " | SAM | OOP" right after " | XML" (end of the skills line).
$r = $d.Content
$r.Find.Execute(" | XML", $true, $false, $false, $false, $false, `
                $true, 1, $false, "", 0)
$ins = $d.Range($r.End, $r.End)
$ins.InsertAfter(" | SAM | OOP")

# 3. Shorten "Search Engine Optimization (SEO)" to "SEO".
$d.Content.Find.Execute("Search Engine Optimization (SEO)", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "SEO", 2)

# 4. Merge the SEO/Progressive-Web-Apps paragraph into the next one by
#    deleting the paragraph mark that separates them.
$r = $d.Content
$r.Find.Execute("Progressive Web Apps", $true, $false, $false, $false, $false, `
                $true, 1, $false, "", 0)
$pm = $d.Range($r.End, $r.End + 1)
$pm.Delete()

# 5. Shorten "User Experience (UX)" to " | Agile | Scrum | UX".
$d.Content.Find.Execute("User Experience (UX)", $true, $false, $false, $false, $false, `
                         $true, 1, $false, " | Agile | Scrum | UX", 2)

# 6. Shorten "User Interface (UI) " to "UI ".
$d.Content.Find.Execute("User Interface (UI) ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "UI ", 2)
